$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "F2"
$ws.Cells.Item(2,3).Value = "F2rl2"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.654227
$ws.Cells.Item(2,8).Value = 4.962681
$ws.Cells.Item(2,9).Value = 0.4107585939979205
$ws.Cells.Item(2,10).Value = 0.4107585939979205
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.390102
$ws.Cells.Item(2,14).Value = 1.170306
$ws.Cells.Item(2,15).Value = 0.09609693103136531
$ws.Cells.Item(2,16).Value = 0.09609693103136531
$ws.Cells.Item(2,17).Value = 0.645317261154
$ws.Cells.Item(2,18).Value = 5.807855350386
$ws.Cells.Item(2,19).Value = 0.03947264027795876
$ws.Cells.Item(2,20).Value = 0.03947264027795876

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "F2"
$ws.Cells.Item(3,3).Value = "F2rl2"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.654227
$ws.Cells.Item(3,8).Value = 4.962681
$ws.Cells.Item(3,9).Value = 0.4107585939979205
$ws.Cells.Item(3,10).Value = 0.4107585939979205
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.860424666666667
$ws.Cells.Item(3,14).Value = 5.581274
$ws.Cells.Item(3,15).Value = 0.4582932178807528
$ws.Cells.Item(3,16).Value = 0.4582932178807528
$ws.Cells.Item(3,17).Value = 3.077564715066
$ws.Cells.Item(3,18).Value = 27.698082435594
$ws.Cells.Item(3,19).Value = 0.1882478778154807
$ws.Cells.Item(3,20).Value = 0.1882478778154807

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "F2"
$ws.Cells.Item(4,3).Value = "F2rl2"
$ws.Cells.Item(4,4).Value = "Resolving-Mac"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.654227
$ws.Cells.Item(4,8).Value = 4.962681
$ws.Cells.Item(4,9).Value = 0.4107585939979205
$ws.Cells.Item(4,10).Value = 0.4107585939979205
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.808937
$ws.Cells.Item(4,14).Value = 5.426811000000001
$ws.Cells.Item(4,15).Value = 0.4456098510878819
$ws.Cells.Item(4,16).Value = 0.4456098510878818
$ws.Cells.Item(4,17).Value = 2.992392426699
$ws.Cells.Item(4,18).Value = 26.931531840291
$ws.Cells.Item(4,19).Value = 0.1830380759044811
$ws.Cells.Item(4,20).Value = 0.1830380759044811

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "F2"
$ws.Cells.Item(5,3).Value = "F2rl2"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.288726
$ws.Cells.Item(5,8).Value = 3.866178
$ws.Cells.Item(5,9).Value = 0.3200015957958394
$ws.Cells.Item(5,10).Value = 0.3200015957958394
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.390102
$ws.Cells.Item(5,14).Value = 1.170306
$ws.Cells.Item(5,15).Value = 0.09609693103136531
$ws.Cells.Item(5,16).Value = 0.09609693103136531
$ws.Cells.Item(5,17).Value = 0.5027345900519999
$ws.Cells.Item(5,18).Value = 4.524611310468
$ws.Cells.Item(5,19).Value = 0.03075117128111962
$ws.Cells.Item(5,20).Value = 0.03075117128111962

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "F2"
$ws.Cells.Item(6,3).Value = "F2rl2"
$ws.Cells.Item(6,4).Value = "MuSCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.288726
$ws.Cells.Item(6,8).Value = 3.866178
$ws.Cells.Item(6,9).Value = 0.3200015957958394
$ws.Cells.Item(6,10).Value = 0.3200015957958394
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.860424666666667
$ws.Cells.Item(6,14).Value = 5.581274
$ws.Cells.Item(6,15).Value = 0.4582932178807528
$ws.Cells.Item(6,16).Value = 0.4582932178807528
$ws.Cells.Item(6,17).Value = 2.397577638974666
$ws.Cells.Item(6,18).Value = 21.578198750772
$ws.Cells.Item(6,19).Value = 0.1466545610642512
$ws.Cells.Item(6,20).Value = 0.1466545610642512

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "F2"
$ws.Cells.Item(7,3).Value = "F2rl2"
$ws.Cells.Item(7,4).Value = "Resolving-Mac"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.288726
$ws.Cells.Item(7,8).Value = 3.866178
$ws.Cells.Item(7,9).Value = 0.3200015957958394
$ws.Cells.Item(7,10).Value = 0.3200015957958394
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.808937
$ws.Cells.Item(7,14).Value = 5.426811000000001
$ws.Cells.Item(7,15).Value = 0.4456098510878819
$ws.Cells.Item(7,16).Value = 0.4456098510878818
$ws.Cells.Item(7,17).Value = 2.331224144262
$ws.Cells.Item(7,18).Value = 20.981017298358
$ws.Cells.Item(7,19).Value = 0.1425958634504686
$ws.Cells.Item(7,20).Value = 0.1425958634504685

$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "F2"
$ws.Cells.Item(8,3).Value = "F2rl2"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.8858993333333333
$ws.Cells.Item(8,8).Value = 2.657698
$ws.Cells.Item(8,9).Value = 0.2199763179924491
$ws.Cells.Item(8,10).Value = 0.2199763179924491
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.390102
$ws.Cells.Item(8,14).Value = 1.170306
$ws.Cells.Item(8,15).Value = 0.09609693103136531
$ws.Cells.Item(8,16).Value = 0.09609693103136531
$ws.Cells.Item(8,17).Value = 0.345591101732
$ws.Cells.Item(8,18).Value = 3.110319915588
$ws.Cells.Item(8,19).Value = 0.02113904905865407
$ws.Cells.Item(8,20).Value = 0.02113904905865407

$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "F2"
$ws.Cells.Item(9,3).Value = "F2rl2"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.8858993333333333
$ws.Cells.Item(9,8).Value = 2.657698
$ws.Cells.Item(9,9).Value = 0.2199763179924491
$ws.Cells.Item(9,10).Value = 0.2199763179924491
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.860424666666667
$ws.Cells.Item(9,14).Value = 5.581274
$ws.Cells.Item(9,15).Value = 0.4582932178807528
$ws.Cells.Item(9,16).Value = 0.4582932178807528
$ws.Cells.Item(9,17).Value = 1.648148971916889
$ws.Cells.Item(9,18).Value = 14.833340747252
$ws.Cells.Item(9,19).Value = 0.1008136546303193
$ws.Cells.Item(9,20).Value = 0.1008136546303192

$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "F2"
$ws.Cells.Item(10,3).Value = "F2rl2"
$ws.Cells.Item(10,4).Value = "Resolving-Mac"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.8858993333333333
$ws.Cells.Item(10,8).Value = 2.657698
$ws.Cells.Item(10,9).Value = 0.2199763179924491
$ws.Cells.Item(10,10).Value = 0.2199763179924491
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.808937
$ws.Cells.Item(10,14).Value = 5.426811000000001
$ws.Cells.Item(10,15).Value = 0.4456098510878819
$ws.Cells.Item(10,16).Value = 0.4456098510878818
$ws.Cells.Item(10,17).Value = 1.602536082342
$ws.Cells.Item(10,18).Value = 14.422824741078
$ws.Cells.Item(10,19).Value = 0.09802361430347581
$ws.Cells.Item(10,20).Value = 0.0980236143034758

$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "F2"
$ws.Cells.Item(11,3).Value = "F2rl2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.1983963333333333
$ws.Cells.Item(11,8).Value = 0.595189
$ws.Cells.Item(11,9).Value = 0.04926349221379096
$ws.Cells.Item(11,10).Value = 0.04926349221379096
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.390102
$ws.Cells.Item(11,14).Value = 1.170306
$ws.Cells.Item(11,15).Value = 0.09609693103136531
$ws.Cells.Item(11,16).Value = 0.09609693103136531
$ws.Cells.Item(11,17).Value = 0.077394806426
$ws.Cells.Item(11,18).Value = 0.696553257834
$ws.Cells.Item(11,19).Value = 0.004734070413632872
$ws.Cells.Item(11,20).Value = 0.004734070413632872

$ws.Cells.Item(12,1).Value = "Resolving-Mac"
$ws.Cells.Item(12,2).Value = "F2"
$ws.Cells.Item(12,3).Value = "F2rl2"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 0.6666666666666666
$ws.Cells.Item(12,7).Value = 0.1983963333333333
$ws.Cells.Item(12,8).Value = 0.595189
$ws.Cells.Item(12,9).Value = 0.04926349221379096
$ws.Cells.Item(12,10).Value = 0.04926349221379096
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 1.860424666666667
$ws.Cells.Item(12,14).Value = 5.581274
$ws.Cells.Item(12,15).Value = 0.4582932178807528
$ws.Cells.Item(12,16).Value = 0.4582932178807528
$ws.Cells.Item(12,17).Value = 0.3691014323095555
$ws.Cells.Item(12,18).Value = 3.321912890785999
$ws.Cells.Item(12,19).Value = 0.02257712437070167
$ws.Cells.Item(12,20).Value = 0.02257712437070167

$ws.Cells.Item(13,1).Value = "Resolving-Mac"
$ws.Cells.Item(13,2).Value = "F2"
$ws.Cells.Item(13,3).Value = "F2rl2"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 0.6666666666666666
$ws.Cells.Item(13,7).Value = 0.1983963333333333
$ws.Cells.Item(13,8).Value = 0.595189
$ws.Cells.Item(13,9).Value = 0.04926349221379096
$ws.Cells.Item(13,10).Value = 0.04926349221379096
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.808937
$ws.Cells.Item(13,14).Value = 5.426811000000001
$ws.Cells.Item(13,15).Value = 0.4456098510878819
$ws.Cells.Item(13,16).Value = 0.4456098510878818
$ws.Cells.Item(13,17).Value = 0.358886468031
$ws.Cells.Item(13,18).Value = 3.229978212279
$ws.Cells.Item(13,19).Value = 0.02195229742945642
$ws.Cells.Item(13,20).Value = 0.021952297429456417

# remove now-obsolete rows 14-17 (Resolving-Mac x Resolving-Mac block superseded)
$ws.Rows("14:17").Delete()
